# Auto-generated edit script applying the cryptos.xlsx data refresh described in the diff.
# Updates Price (D) and Volume(1h) (E) columns for all 50 coin rows, plus a couple of
# coin-identity swaps (rows 31/32 PancakeSwap<->Filecoin, rows 48/49 Maker<->Aptos).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Some new Price values are plain decimal numbers (e.g. "235.42"). The source data stores
# prices as literal text (it mixes plain decimals with thousand-separated "." numbers like
# "30.180.34", which can't be a real number), so force those specific cells to Text format
# first - otherwise Excel would silently convert them to numeric values.
$textCells = @("D4", "D5", "D6", "D7", "D8", "D9", "D11", "D12", "D15", "D16", "D18", "D19", "D20", "D22", "D23", "D24", "D25", "D26", "D27", "D28", "D29", "D30", "D31", "D32", "D33", "D34", "D35", "D36", "D39", "D40", "D42", "D43", "D44", "D45", "D46", "D47", "D48", "D49", "D50", "D51")
foreach ($c in $textCells) { $ws.Range($c).NumberFormat = "@" }

# Write the updated values
$ws.Range("D4").Value = "0.9993"
$ws.Range("D5").Value = "235.42"
$ws.Range("D6").Value = "0.9989"
$ws.Range("D7").Value = "0.4691"
$ws.Range("D8").Value = "0.2888"
$ws.Range("D9").Value = "0.06568"
$ws.Range("D11").Value = "0.07977"
$ws.Range("D12").Value = "97.38"
$ws.Range("D15").Value = "0.6769"
$ws.Range("D16").Value = "268.08"
$ws.Range("D18").Value = "13.60"
$ws.Range("D19").Value = "0.000007663"
$ws.Range("D20").Value = "0.9990"
$ws.Range("D22").Value = "0.9998"
$ws.Range("D23").Value = "5.205"
$ws.Range("D24").Value = "6.154"
$ws.Range("D25").Value = "166.78"
$ws.Range("D26").Value = "9.179"
$ws.Range("D27").Value = "18.92"
$ws.Range("D28").Value = "1.940"
$ws.Range("D29").Value = "1.382"
$ws.Range("D30").Value = "0.09876"
$ws.Range("D31").Value = "4.317"
$ws.Range("D32").Value = "1.466"
$ws.Range("D33").Value = "4.026"
$ws.Range("D34").Value = "0.04708"
$ws.Range("D35").Value = "1.121"
$ws.Range("D36").Value = "0.7001"
$ws.Range("D39").Value = "2.600"
$ws.Range("D40").Value = "6.326"
$ws.Range("D42").Value = "1.931"
$ws.Range("D43").Value = "0.8394"
$ws.Range("D44").Value = "0.9986"
$ws.Range("D45").Value = "103.53"
$ws.Range("D46").Value = "0.4141"
$ws.Range("D47").Value = "9.170"
$ws.Range("D48").Value = "7.041"
$ws.Range("D49").Value = "935.65"
$ws.Range("D50").Value = "33.97"
$ws.Range("D51").Value = "0.05657"
$ws.Range("D2").Value = "30.180.34"
$ws.Range("E2").Value = "  -0.50%  "
$ws.Range("D3").Value = "1.857.41"
$ws.Range("E4").Value = "  -0.16%  "
$ws.Range("E5").Value = "  -0.11%  "
$ws.Range("E7").Value = "  +0.22%  "
$ws.Range("E8").Value = "  +1.47%  "
$ws.Range("E9").Value = "  +0.36%  "
$ws.Range("E10").Value = "  +1.86%  "
$ws.Range("E11").Value = "  +1.19%  "
$ws.Range("E12").Value = "  -0.51%  "
$ws.Range("D13").Value = "1.857.10"
$ws.Range("E13").Value = "  -0.64%  "
$ws.Range("E14").Value = "  +0.09%  "
$ws.Range("E15").Value = "  +0.14%  "
$ws.Range("E16").Value = "  -3.15%  "
$ws.Range("D17").Value = "30.146.29"
$ws.Range("E17").Value = "  -0.59%  "
$ws.Range("E18").Value = "  +6.88%  "
$ws.Range("E19").Value = "  +4.87%  "
$ws.Range("E20").Value = "  -0.19%  "
$ws.Range("D21").Value = "2.098.95"
$ws.Range("E21").Value = "  -0.86%  "
$ws.Range("E22").Value = "  -0.12%  "
$ws.Range("E23").Value = "  -4.85%  "
$ws.Range("E24").Value = "  +0.08%  "
$ws.Range("E25").Value = "  +0.71%  "
$ws.Range("E26").Value = "  +0.41%  "
$ws.Range("E27").Value = "  -0.95%  "
$ws.Range("E28").Value = "  +0.54%  "
$ws.Range("E29").Value = "  +0.12%  "
$ws.Range("E30").Value = "  +2.67%  "
$ws.Range("B31").Value = "Filecoin"
$ws.Range("C31").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("E31").Value = "  -1.59%  "
$ws.Range("B32").Value = "PancakeSwap"
$ws.Range("C32").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("E32").Value = "  -0.72%  "
$ws.Range("E33").Value = "  -1.62%  "
$ws.Range("E34").Value = "  +0.20%  "
$ws.Range("E35").Value = "  -0.59%  "
$ws.Range("E36").Value = "  -0.75%  "
$ws.Range("E37").Value = "  -0.79%  "
$ws.Range("E38").Value = "  +0.72%  "
$ws.Range("E39").Value = "  +2.55%  "
$ws.Range("E40").Value = "  +0.08%  "
$ws.Range("E41").Value = "  -0.81%  "
$ws.Range("E42").Value = "  -1.15%  "
$ws.Range("E43").Value = "  -1.19%  "
$ws.Range("E44").Value = "  -0.20%  "
$ws.Range("E45").Value = "  -0.27%  "
$ws.Range("E46").Value = "  -1.09%  "
$ws.Range("E47").Value = "  -0.98%  "
$ws.Range("B48").Value = "Aptos"
$ws.Range("C48").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("E48").Value = "  -2.14%  "
$ws.Range("B49").Value = "Maker"
$ws.Range("C49").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("E49").Value = "  -0.03%  "
$ws.Range("E50").Value = "  -0.61%  "
$ws.Range("E51").Value = "  +0.41%  "
